# vividus-plugin-excel TestTemplate.xlsx fixture update (#5166)
# "Fix handling of cells without content":
# adds a new "Null" column (header only, no data rows) to the
# DifferentTypes sheet, and leaves that sheet active/selected with the
# cursor parked just past the used range.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("DifferentTypes")

# New header-only column E ("Null") - rows 2/3 stay empty to exercise
# handling of cells without content.
$ws.Range("E1").Value = "Null"

# Make this sheet the active tab/selection, cursor below the new column.
$ws.Activate() | Out-Null
$ws.Range("E6").Select() | Out-Null
